# Natmi following Dr Hou advice:
# Recompute the Sema4g-Plxnb2 LR-pairs table so the sending/target clusters
# sweep all four cell types (ECs, FAPs, M2, sCs) instead of just three,
# and refresh every derived statistic (detection rate, expression,
# specificity and edge-weight columns) accordingly. This grows the table
# from 12 data rows (3x4) to 16 data rows (4x4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
  @("ECs", "Sema4g", "Plxnb2", "ECs", 1, 0.3333333333333333, 0.5158756666666667, 1.547627, 0.08115675834014567, 0.08115675834014567, 3, 1, 10.39091433333333, 31.172743, 0.0835098648954196, 0.0835098648954196, 5.360419858984556, 48.243778730861, 0.006777389924335783, 0.006777389924335783),
  @("ECs", "Sema4g", "Plxnb2", "FAPs", 1, 0.3333333333333333, 0.5158756666666667, 1.547627, 0.08115675834014567, 0.08115675834014567, 3, 1, 35.71561933333334, 107.146858, 0.287039855156433, 0.287039855156433, 18.42481893399622, 165.823370405966, 0.02329522415892105, 0.02329522415892105),
  @("ECs", "Sema4g", "Plxnb2", "M2", 1, 0.3333333333333333, 0.5158756666666667, 1.547627, 0.08115675834014567, 0.08115675834014567, 3, 1, 51.87044599999999, 155.611338, 0.4168732219867682, 0.4168732219867682, 26.75870091054733, 240.828308194926, 0.03383207933525804, 0.03383207933525805),
  @("ECs", "Sema4g", "Plxnb2", "sCs", 1, 0.3333333333333333, 0.5158756666666667, 1.547627, 0.08115675834014567, 0.08115675834014567, 3, 1, 26.450408, 79.351224, 0.2125770579613792, 0.2125770579613792, 13.64512186060533, 122.806096745448, 0.01725206492163079, 0.01725206492163079),
  @("FAPs", "Sema4g", "Plxnb2", "ECs", 3, 1, 1.828868333333333, 5.486604999999999, 0.2877147246027853, 0.2877147246027853, 3, 1, 10.39091433333333, 31.172743, 0.0835098648954196, 0.0835098648954196, 19.00361417861278, 171.032527607515, 0.02402701778000146, 0.02402701778000146),
  @("FAPs", "Sema4g", "Plxnb2", "FAPs", 3, 1, 1.828868333333333, 5.486604999999999, 0.2877147246027853, 0.2877147246027853, 3, 1, 35.71561933333334, 107.146858, 0.287039855156433, 0.287039855156433, 65.3191652041211, 587.8724868370899, 0.0825855928763565, 0.0825855928763565),
  @("FAPs", "Sema4g", "Plxnb2", "M2", 3, 1, 1.828868333333333, 5.486604999999999, 0.2877147246027853, 0.2877147246027853, 3, 1, 51.87044599999999, 155.611338, 0.4168732219867682, 0.4168732219867682, 94.86421612527664, 853.7779451274898, 0.1199405642581988, 0.1199405642581988),
  @("FAPs", "Sema4g", "Plxnb2", "sCs", 3, 1, 1.828868333333333, 5.486604999999999, 0.2877147246027853, 0.2877147246027853, 3, 1, 26.450408, 79.351224, 0.2125770579613792, 0.2125770579613792, 48.37431359494666, 435.3688223545199, 0.06116154968822855, 0.06116154968822855),
  @("M2", "Sema4g", "Plxnb2", "ECs", 3, 1, 0.2256573333333333, 0.6769719999999999, 0.03550006106577688, 0.03550006106577688, 3, 1, 10.39091433333333, 31.172743, 0.0835098648954196, 0.0835098648954196, 2.34478601935511, 21.103074174196, 0.002964605303382173, 0.002964605303382173),
  @("M2", "Sema4g", "Plxnb2", "FAPs", 3, 1, 0.2256573333333333, 0.6769719999999999, 0.03550006106577688, 0.03550006106577688, 3, 1, 35.71561933333334, 107.146858, 0.287039855156433, 0.287039855156433, 8.059491417108443, 72.535422753976, 0.01018993238636512, 0.01018993238636512),
  @("M2", "Sema4g", "Plxnb2", "M2", 3, 1, 0.2256573333333333, 0.6769719999999999, 0.03550006106577688, 0.03550006106577688, 3, 1, 51.87044599999999, 155.611338, 0.4168732219867682, 0.4168732219867682, 11.70494652317066, 105.344518708536, 0.01479902483721743, 0.01479902483721744),
  @("M2", "Sema4g", "Plxnb2", "sCs", 3, 1, 0.2256573333333333, 0.6769719999999999, 0.03550006106577688, 0.03550006106577688, 3, 1, 26.450408, 79.351224, 0.2125770579613792, 0.2125770579613792, 5.968728534858665, 53.71855681372799, 0.007546498538812154, 0.007546498538812155),
  @("sCs", "Sema4g", "Plxnb2", "ECs", 3, 1, 3.786132333333333, 11.358397, 0.5956284559912921, 0.5956284559912921, 3, 1, 10.39091433333333, 31.172743, 0.0835098648954196, 0.0835098648954196, 39.34137673033011, 354.072390572971, 0.04974085188770019, 0.04974085188770019),
  @("sCs", "Sema4g", "Plxnb2", "FAPs", 3, 1, 3.786132333333333, 11.358397, 0.5956284559912921, 0.5956284559912921, 3, 1, 35.71561933333334, 107.146858, 0.287039855156433, 0.287039855156433, 135.2240611629585, 1217.016550466626, 0.1709691057347903, 0.1709691057347903),
  @("sCs", "Sema4g", "Plxnb2", "M2", 3, 1, 3.786132333333333, 11.358397, 0.5956284559912921, 0.5956284559912921, 3, 1, 51.87044599999999, 155.611338, 0.4168732219867682, 0.4168732219867682, 196.3883727450206, 1767.495354705186, 0.2483015535560939, 0.2483015535560939),
  @("sCs", "Sema4g", "Plxnb2", "sCs", 3, 1, 3.786132333333333, 11.358397, 0.5956284559912921, 0.5956284559912921, 3, 1, 26.450408, 79.351224, 0.2125770579613792, 0.2125770579613792, 100.1447449586587, 901.302704627928, 0.1266169448127077, 0.1266169448127077)
)

$r = 2
foreach ($row in $rowsData) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value2 = $val
        $c++
    }
    $r++
}
